$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at C (shifts all_audio_played..worker_id one column right)
# and populate it with the new "accepted_and_use" header.
$ws.Columns.Item(3).Insert()
$ws.Range("C1").Value = "accepted_and_use"

# Insert a new column at G (old correct_math's new slot), shifting
# correct_math/correct_tps/worker_id one column right, and populate it
# with "correct_gold_question".
$ws.Columns.Item(7).Insert()
$ws.Range("G1").Value = "correct_gold_question"

# Insert a new column at J (worker_id's current slot), shifting worker_id
# one column right, and populate it with "variance_in_ratings".
$ws.Columns.Item(10).Insert()
$ws.Range("J1").Value = "variance_in_ratings"

# Remove the old sample data row; only the header row remains.
$ws.Rows.Item(2).Delete()
